$wb = $excel.ActiveWorkbook

# Add a new worksheet and move it to the first position, named "Player Info"
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($wb.Worksheets.Item(1))

# Headers
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header styling used on the other sheets
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row (force text storage so numeric-looking strings stay strings,
# then restore the default "Normal" style so no stray formatting is left
# behind on the cell)
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6664"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "David Alan Payne"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# Update ODI Batting sheet: rename MATCH_CARD_LINK -> MATCH_CODE, update value
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4602"
$batting.Range("D2").Style = "Normal"

# Update ODI Bowling sheet: rename MATCH_CARD_LINK -> MATCH_CODE, update value
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4602"
$bowling.Range("B2").Style = "Normal"
